$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F (6) through V (22) hold the match data that gets reshuffled
# between rows when the source data is re-scraped/re-sorted; columns
# A-E (Indice/pais/torneio/temporada/data_partida) stay put per row.
$firstCol = 6
$lastCol = 22

function Swap-Rows($rowA, $rowB) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $valA = $ws.Cells.Item($rowA, $c).Value2
        $valB = $ws.Cells.Item($rowB, $c).Value2
        $ws.Cells.Item($rowA, $c).Value2 = $valB
        $ws.Cells.Item($rowB, $c).Value2 = $valA
    }
}

# Simple pairwise swaps
Swap-Rows 2 3
Swap-Rows 6 7
Swap-Rows 18 19
Swap-Rows 28 29
Swap-Rows 30 31
Swap-Rows 48 49

# Rows 53-55 rotate up by one: old 54 -> 53, old 55 -> 54, old 53 -> 55
$tmp53 = @{}
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $tmp53[$c] = $ws.Cells.Item(53, $c).Value2
}
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item(53, $c).Value2 = $ws.Cells.Item(54, $c).Value2
}
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item(54, $c).Value2 = $ws.Cells.Item(55, $c).Value2
}
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws.Cells.Item(55, $c).Value2 = $tmp53[$c]
}

# Append the new match row (row 132) at the end of the sheet.
# Copy formatting from the last existing row (131) first so the new
# row's styled cells (A: bold/bordered index style, E: date style)
# reuse the same style indices as the rest of the sheet.
$newRow = 132
$prevRow = $newRow - 1
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($prevRow, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value2 = 131
$ws.Cells.Item($newRow, 2).Value2 = "italy"
$ws.Cells.Item($newRow, 3).Value2 = "serie-a"
$ws.Cells.Item($newRow, 4).Value2 = "2023-2024"
$ws.Cells.Item($newRow, 5).Value2 = 45261.86458333334
$ws.Cells.Item($newRow, 6).Value2 = "Monza"
$ws.Cells.Item($newRow, 7).Value2 = 1
$ws.Cells.Item($newRow, 8).Value2 = "Juventus"
$ws.Cells.Item($newRow, 9).Value2 = 2
$ws.Cells.Item($newRow, 10).Value2 = 3.34
$ws.Cells.Item($newRow, 11).Value2 = "15/11/2023 16:01"
$ws.Cells.Item($newRow, 12).Value2 = 4.37
$ws.Cells.Item($newRow, 13).Value2 = "01/12/2023 20:44"
$ws.Cells.Item($newRow, 14).Value2 = 3.38
$ws.Cells.Item($newRow, 15).Value2 = "15/11/2023 16:01"
$ws.Cells.Item($newRow, 16).Value2 = 3.43
$ws.Cells.Item($newRow, 17).Value2 = "01/12/2023 20:44"
$ws.Cells.Item($newRow, 18).Value2 = 2.25
$ws.Cells.Item($newRow, 19).Value2 = "15/11/2023 16:01"
$ws.Cells.Item($newRow, 20).Value2 = 1.93
$ws.Cells.Item($newRow, 21).Value2 = "01/12/2023 20:07"
$ws.Cells.Item($newRow, 22).Value2 = "https://www.betexplorer.com/football/italy/serie-a/monza-juventus/EXHeCNVs/"
